# Insert a new weekly price record at row 94 ("Hortaliza, Vega Monumental
# Concepción - Zapallo"). All existing records from row 94 down to row 159
# shift down by one row (to rows 95-160), and the new record is written
# into the freed-up row 94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 94..159 down to 95..160, leaving a blank (but style-carrying) row 94
$ws.Rows.Item(94).Insert()

# Populate the new row 94 with the new data record
$ws.Range("A94").Value = 11
$ws.Range("B94").Value = 'Vega Monumental Concepción'
$ws.Range("C94").Value = 'Bíobío'
$ws.Range("D94").Value = 44572
$ws.Range("E94").Value = 8
$ws.Range("F94").Value = 100112045
$ws.Range("G94").Value = 'Zapallo'
$ws.Range("H94").Value = 'Camote'
$ws.Range("I94").Value = '1a nueva(o)'
$ws.Range("J94").Value = 290
$ws.Range("K94").Value = 300
$ws.Range("L94").Value = 350
$ws.Range("M94").Value = 326
$ws.Range("N94").Value = '$/kilo (volumen en unidades)'
$ws.Range("O94").Value = 'Región Metropolitana'
$ws.Range("P94").Value = 326
$ws.Range("Q94").Value = 1
$ws.Range("R94").Value = 'Hortaliza'
